# Update countries & provincias Spain
#
# The "Pais" table (A3:H216) is kept sorted descending by column B
# ("Casos totales"). Three countries (Rusia, Banglades, Guyana) receive
# refreshed case numbers that move them further up the ranking. Rather
# than inserting new rows, the existing rows are reused: starting at the
# row each country now belongs in, the country name (column A) and data
# (columns B:H) of every row below shift down by one, until reaching a
# row that already held the right data; the freed-up row gets the
# country's fresh figures. A few unrelated rows (Moldavia, Estonia,
# Haiti) simply get corrected case counts in place, with no row shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rusia moves up into row 14; Brasil/Canada/Paises Bajos shift down one row ---
$ws.Range("A14").Value = "Rusia"
$ws.Range("B14").Value = 32008
$ws.Range("C14").Value = 4070
$ws.Range("D14").Value = 2590
$ws.Range("E14").Value = 29145
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 41
$ws.Range("H14").Value = 273

$ws.Range("A15").Value = "Brasil"
$ws.Range("B15").Value = 30891
$ws.Range("C15").Value = 208
$ws.Range("D15").Value = 14026
$ws.Range("E15").Value = 14913
$ws.Range("F15").Value = 6634
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 1952

$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 30106
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 9729
$ws.Range("E16").Value = 19182
$ws.Range("F16").Value = 557
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 1195

$ws.Range("A17").Value = "Paises Bajos"
$ws.Range("B17").Value = 29214
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 250
$ws.Range("E17").Value = 25649
$ws.Range("F17").Value = 1279
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 3315

# --- Moldavia: corrected active/recovered counts (no row shift) ---
$ws.Range("D60").Value = 276
$ws.Range("E60").Value = 1824

# --- Banglades moves up into row 61; Croacia/Hungria/Islandia/Barein shift down one row ---
$ws.Range("A61").Value = "Banglades"
$ws.Range("B61").Value = 1838
$ws.Range("C61").Value = 266
$ws.Range("D61").Value = 58
$ws.Range("E61").Value = 1705
$ws.Range("F61").Value = 1
$ws.Range("G61").Value = 15
$ws.Range("H61").Value = 75

$ws.Range("A62").Value = "Croacia"
$ws.Range("B62").Value = 1791
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 529
$ws.Range("E62").Value = 1227
$ws.Range("F62").Value = 31
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 35

$ws.Range("A63").Value = "Hungria"
$ws.Range("B63").Value = 1763
$ws.Range("C63").Value = 111
$ws.Range("D63").Value = 207
$ws.Range("E63").Value = 1400
$ws.Range("F63").Value = 63
$ws.Range("G63").Value = 14
$ws.Range("H63").Value = 156

$ws.Range("A64").Value = "Islandia"
$ws.Range("B64").Value = 1739
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 1144
$ws.Range("E64").Value = 587
$ws.Range("F64").Value = 6
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 8

$ws.Range("A65").Value = "Barein"
$ws.Range("B65").Value = 1700
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 703
$ws.Range("E65").Value = 990
$ws.Range("F65").Value = 3
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 7

# --- Estonia: corrected case counts (no row shift) ---
$ws.Range("A68").Value = "Estonia"
$ws.Range("B68").Value = 1459
$ws.Range("C68").Value = 25
$ws.Range("D68").Value = 145
$ws.Range("E68").Value = 1276
$ws.Range("F68").Value = 11
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 38

# --- Guyana moves up into row 150; San Martin/Cabo Verde/Polinesia shift down one row ---
$ws.Range("A150").Value = "Guyana"
$ws.Range("B150").Value = 57
$ws.Range("C150").Value = 2
$ws.Range("D150").Value = 9
$ws.Range("E150").Value = 42
$ws.Range("F150").Value = 4
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 6

$ws.Range("A151").Value = "San Martin (Parte Holandesa)"
$ws.Range("B151").Value = 57
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 12
$ws.Range("E151").Value = 36
$ws.Range("F151").Value = 6
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 9

$ws.Range("A152").Value = "Cabo Verde"
$ws.Range("B152").Value = 56
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 1
$ws.Range("E152").Value = 54
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 1

$ws.Range("A153").Value = "Polinesia Francesa"
$ws.Range("B153").Value = 55
$ws.Range("C153").Value = 0
$ws.Range("D153").Value = 0
$ws.Range("E153").Value = 55
$ws.Range("F153").Value = 1
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 0

# --- Haiti: corrected case counts (no row shift) ---
$ws.Range("B161").Value = 43
$ws.Range("C161").Value = 2
$ws.Range("D161").Value = 40
$ws.Range("E161").Value = 38
